$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text format on price cells whose new values would otherwise be
# auto-detected as numbers (the column stores prices as literal text).
$textForceCells = @('D5', 'D6', 'D11', 'D14', 'D20', 'D21', 'D27', 'D31', 'D40', 'D42', 'D47')
foreach ($addr in $textForceCells) {
    $ws.Range($addr).NumberFormat = "@"
}

# Apply the updated Price (D) and Volume(1h) (E) values row by row.
$ws.Range("D2").Value2 = '63.727.59'
$ws.Range("E2").Value2 = '  +3.12%  '
$ws.Range("D3").Value2 = '3.126.62'
$ws.Range("E3").Value2 = '  +1.70%  '
$ws.Range("E4").Value2 = '  -0.01%  '
$ws.Range("D5").Value2 = '587.98'
$ws.Range("E5").Value2 = '  +1.51%  '
$ws.Range("D6").Value2 = '147.15'
$ws.Range("E6").Value2 = '  +3.81%  '
$ws.Range("E7").Value2 = '  +0.01%  '
$ws.Range("D8").Value2 = '3.121.24'
$ws.Range("E8").Value2 = '  +1.85%  '
$ws.Range("E9").Value2 = '  +0.91%  '
$ws.Range("E10").Value2 = '  +15.15%  '
$ws.Range("D11").Value2 = '5.69'
$ws.Range("E11").Value2 = '  +0.50%  '
$ws.Range("E12").Value2 = '  +0.96%  '
$ws.Range("E13").Value2 = '  +5.41%  '
$ws.Range("D14").Value2 = '37.12'
$ws.Range("E14").Value2 = '  +5.66%  '
$ws.Range("E15").Value2 = '  -0.64%  '
$ws.Range("D16").Value2 = '3.644.58'
$ws.Range("E16").Value2 = '  +1.80%  '
$ws.Range("E17").Value2 = '  -0.94%  '
$ws.Range("D18").Value2 = '63.606.09'
$ws.Range("E18").Value2 = '  +3.05%  '
$ws.Range("D19").Value2 = '3.122.04'
$ws.Range("E19").Value2 = '  +1.74%  '
$ws.Range("D20").Value2 = '465.01'
$ws.Range("E20").Value2 = '  +4.06%  '
$ws.Range("D21").Value2 = '14.36'
$ws.Range("E21").Value2 = '  +2.71%  '
$ws.Range("E22").Value2 = '  -0.01%  '
$ws.Range("E23").Value2 = '  +1.46%  '
$ws.Range("E24").Value2 = '  -3.06%  '
$ws.Range("E25").Value2 = '  +0.88%  '
$ws.Range("E26").Value2 = '  +0.01%  '
$ws.Range("D27").Value2 = '9.02'
$ws.Range("E27").Value2 = '  +10.20%  '
$ws.Range("E28").Value2 = '  +1.84%  '
$ws.Range("E29").Value2 = '  -1.63%  '
$ws.Range("E30").Value2 = '  -0.05%  '
$ws.Range("D31").Value2 = '6.87'
$ws.Range("E31").Value2 = '  +1.47%  '
$ws.Range("E32").Value2 = '  +1.47%  '
$ws.Range("E33").Value2 = '  -4.52%  '
$ws.Range("D34").Value2 = '0.0₃0879'
$ws.Range("E34").Value2 = '  +10.94%  '
$ws.Range("E35").Value2 = '  +8.67%  '
$ws.Range("E36").Value2 = '  +1.90%  '
$ws.Range("E37").Value2 = '  +16.65%  '
$ws.Range("E38").Value2 = '  +1.57%  '
$ws.Range("E39").Value2 = '  +1.67%  '
$ws.Range("D40").Value2 = '456.68'
$ws.Range("E40").Value2 = '  +8.85%  '
$ws.Range("E41").Value2 = '  -0.70%  '
$ws.Range("D42").Value2 = '0.0374'
$ws.Range("E42").Value2 = '  +1.47%  '
$ws.Range("D43").Value2 = '2.908.49'
$ws.Range("E43").Value2 = '  -1.70%  '
$ws.Range("E44").Value2 = '  +0.78%  '
$ws.Range("E46").Value2 = '  +3.47%  '
$ws.Range("D47").Value2 = '127.35'
$ws.Range("E47").Value2 = '  +2.49%  '
$ws.Range("E48").Value2 = '  +2.42%  '
$ws.Range("E50").Value2 = '  +0.56%  '
$ws.Range("E51").Value2 = '  +1.67%  '
